# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml -> "Integral"      (used by the slide master, i.e. the slides)
#   ppt/theme/theme2.xml -> "Office Theme"  (used by the notes master)
#
# The authored edit swaps the two themes' content wholesale: theme1.xml ends up
# holding the "Office Theme" color scheme and theme2.xml ends up holding the
# "Integral" color scheme (font/format schemes are identical between the two
# themes already, so only the color scheme actually changes visually).
#
# The PowerPoint object model only exposes the *slide master's* theme for
# mutation (Master.ColorScheme), so we recolor it here to the "Office Theme"
# palette, channel by channel, via RGBColor.RGB (stored as 0xBBGGRR).

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme

$cs.Colors(1).RGB  = 0x000000   # dk1      000000
$cs.Colors(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$cs.Colors(3).RGB  = 0x6A5444   # dk2      44546A
$cs.Colors(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$cs.Colors(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$cs.Colors(6).RGB  = 0x317DED   # accent2  ED7D31
$cs.Colors(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$cs.Colors(8).RGB  = 0x00C0FF   # accent4  FFC000
$cs.Colors(9).RGB  = 0xC47244   # accent5  4472C4
$cs.Colors(10).RGB = 0x47AD70   # accent6  70AD47
$cs.Colors(11).RGB = 0xC16305   # hlink    0563C1
$cs.Colors(12).RGB = 0x724F95   # folHlink 954F72
